# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets.
# Both sheets share the same row layout / new values, so iterate over
# each sheet by name and apply the same row -> new-value map.

$updates = @{
    2  = 1148
    3  = 863
    4  = 284
    5  = 52
    6  = 1114
    8  = 2393
    9  = 7780
    10 = 928
    11 = 450
    12 = 386
    13 = 159
    15 = 4
    16 = 163
    17 = 8001
    19 = 1388
    22 = 225
    23 = 178
    24 = 328
    25 = 173
    26 = 167
    28 = 112
    29 = 29
    30 = 427
    31 = 1158
    33 = 98
    35 = 84
}

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
